$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$c = $ws.Range('D2')
$c.NumberFormat = '@'
$c.Value = '64.669.38'
$c.Style = 'Normal'
$ws.Range('E2').Value = '  -1.12%  '

# Row 3
$c = $ws.Range('D3')
$c.NumberFormat = '@'
$c.Value = '3.418.26'
$c.Style = 'Normal'
$ws.Range('E3').Value = '  -2.65%  '

# Row 4
$c = $ws.Range('D4')
$c.NumberFormat = '@'
$c.Value = '0.996'
$c.Style = 'Normal'
$ws.Range('E4').Value = '  -0.39%  '

# Row 5
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '583.06'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  -2.75%  '

# Row 6
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '135.91'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  -4.74%  '

# Row 7
$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '0.997'
$c.Style = 'Normal'
$ws.Range('E7').Value = '  -0.20%  '

# Row 8
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '3.415.28'
$c.Style = 'Normal'
$ws.Range('E8').Value = '  -2.74%  '

# Row 9
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.487'
$c.Style = 'Normal'
$ws.Range('E9').Value = '  -6.19%  '

# Row 10
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '0.119'
$c.Style = 'Normal'
$ws.Range('E10').Value = '  -10.02%  '

# Row 11
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '7.04'
$c.Style = 'Normal'
$ws.Range('E11').Value = '  -9.77%  '

# Row 12
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '0.373'
$c.Style = 'Normal'
$ws.Range('E12').Value = '  -7.30%  '

# Row 13
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '3.980.98'
$c.Style = 'Normal'
$ws.Range('E13').Value = '  -3.20%  '

# Row 14
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '0.0000177'
$c.Style = 'Normal'
$ws.Range('E14').Value = '  -9.21%  '

# Row 15
$ws.Range('B15').Value = 'Avalanche'
$ws.Range('C15').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '26.19'
$c.Style = 'Normal'
$ws.Range('E15').Value = '  -7.38%  '

# Row 16
$ws.Range('B16').Value = 'TRON'
$ws.Range('C16').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '0.115'
$c.Style = 'Normal'
$ws.Range('E16').Value = '  -1.78%  '

# Row 17
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '3.412.80'
$c.Style = 'Normal'
$ws.Range('E17').Value = '  -2.61%  '

# Row 18
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '64.462.33'
$c.Style = 'Normal'
$ws.Range('E18').Value = '  -1.40%  '

# Row 19
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '9.54'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  -11.75%  '

# Row 20
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '5.80'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  -5.56%  '

# Row 21
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '13.47'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  -5.99%  '

# Row 22
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '381.36'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  -8.60%  '

# Row 23
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '0.548'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  -7.58%  '

# Row 24
$ws.Range('E24').Value = '  -0.07%  '

# Row 25
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '71.83'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  -6.64%  '

# Row 26
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '3.542.99'
$c.Style = 'Normal'
$ws.Range('E26').Value = '  -3.07%  '

# Row 27
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '0.0000105'
$c.Style = 'Normal'
$ws.Range('E27').Value = '  -7.64%  '

# Row 28
$ws.Range('E28').Value = '  +0.24%  '

# Row 29
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '7.12'
$c.Style = 'Normal'
$ws.Range('E29').Value = '  -7.64%  '

# Row 30
$ws.Range('B30').Value = 'InternetComputer(DFINITY)'
$ws.Range('C30').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '8.01'
$c.Style = 'Normal'
$ws.Range('E30').Value = '  -9.33%  '

# Row 31
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '2.18'
$c.Style = 'Normal'
$ws.Range('E31').Value = '  -10.72%  '

# Row 32
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '3.416.29'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  -2.93%  '

# Row 33
$ws.Range('E33').Value = '  -0.07%  '

# Row 34
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '22.92'
$c.Style = 'Normal'
$ws.Range('E34').Value = '  -4.95%  '

# Row 35
$ws.Range('E35').Value = '  -9.31%  '

# Row 36
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '167.82'
$c.Style = 'Normal'
$ws.Range('E36').Value = '  -3.40%  '

# Row 37
$ws.Range('B37').Value = 'Fetch.AI'
$ws.Range('C37').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '1.18'
$c.Style = 'Normal'
$ws.Range('E37').Value = '  -10.04%  '

# Row 38
$ws.Range('B38').Value = 'Aptos'
$ws.Range('C38').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '6.71'
$c.Style = 'Normal'
$ws.Range('E38').Value = '  -10.34%  '

# Row 39
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '1.45'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  -7.06%  '

# Row 40
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '4.62'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  -11.23%  '

# Row 41
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '0.0752'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  -6.91%  '

# Row 42
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '0.807'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  -5.30%  '

# Row 43
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '0.996'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  -0.33%  '

# Row 44
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '42.19'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  -6.55%  '

# Row 45
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '4.28'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  -14.13%  '

# Row 46
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '1.60'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  -9.16%  '

# Row 47
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '1.13'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  +3.47%  '

# Row 48
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '22.52'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  -2.84%  '

# Row 49
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '6.42'
$c.Style = 'Normal'
$ws.Range('E49').Value = '  -8.25%  '

# Row 50
$ws.Range('B50').Value = 'dogwifhat'
$ws.Range('C50').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '2.01'
$c.Style = 'Normal'
$ws.Range('E50').Value = '  -14.07%  '

# Row 51
$ws.Range('B51').Value = 'Maker'
$ws.Range('C51').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '2.160.85'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  -7.02%  '
